# Updates cryptocurrency price (D) and 1h-volume-change (E) columns
# to reflect refreshed data, per the "Updated cryptos list ... GitHub Actions" commit.
#
# Price cells (column D) hold plain-text numbers such as "26.351.64" or
# "1.010" (European-style thousands separators / trailing zeros) rather than
# real numerics, so each D write forces the cell to Text format first and
# restores the default "Normal" style afterwards (leaving no residual
# NumberFormat) so the stored value keeps its exact textual form instead of
# being auto-coerced into a number by Excel's usual type inference.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.351.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.671.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5317'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06373'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07853'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.542'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.673.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.899.58'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5626'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.70%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.355.29'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.732'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '202.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.17%  '

$ws.Range("E22").Value = '  +1.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.074'
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1218'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.265'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.516'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.72%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05903'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.02%  '

$ws.Range("E31").Value = '  +0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.537'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.333'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.79%  '

$ws.Range("E34").Value = '  -0.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9678'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.835'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.433'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5813'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01621'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.978'
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.081.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8615'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.010'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.808.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.19%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.62'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈106'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.22%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4419'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.072'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05150'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.31%  '
